$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44357, 976, 3808, 1182, 6585, 438, 2315, 1000, 28500, 0, 0, 0, 3674, 420, 910, 44881),
    @(44358, 1050, 4258, 905, 6485, 496, 2287, 0, 28500, 0, 0, 0, 3674, 61, 959, 45204)
)

$startRow = 111
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
}
